# Update the Sprint 4 date line:
#   "10.1/17.1/24.1/31.1" -> "9.1/16.1/23.1/30.1"
# Word records this kind of multi-spot manual edit by splitting the
# paragraph into several runs around each edited character, so the
# final run layout has 8 separate <w:r> elements instead of 1.

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Text -like "*10.1/17.1/24.1/31.1*") {
        $target = $r
        break
    }
}

$s = $target.Start

# Old text layout (offsets relative to $s):
#   "10"(0,2) ".1/1"(2,4) "7"(6,1) ".1/2"(7,4) "4"(11,1) ".1/3"(12,4) "1"(16,1) ".1"(17,2)
# Replace the three same-length numeric spots first (offsets stay valid
# since their lengths don't change), then the length-changing "10"->"9"
# spot last since it sits left of everything else and would otherwise
# shift every subsequent offset by -1.

$o = $s + 6
$sub = $d.Range($o, $o + 1)
$sub.Text = "6"

$o = $s + 11
$sub = $d.Range($o, $o + 1)
$sub.Text = "3"

$o = $s + 16
$sub = $d.Range($o, $o + 1)
$sub.Text = "0"

$o = $s + 0
$sub = $d.Range($o, $o + 2)
$sub.Text = "9"

# Re-split the paragraph into the 8 runs matching the final text
# "9.1/16.1/23.1/30.1", by round-tripping FormattedText on each
# segment boundary (assigning a range's own FormattedText back onto
# itself forces a clean run split with no leftover rPr).
$lens = @(1, 4, 1, 4, 1, 4, 1, 2)
$pos = $s
foreach ($len in $lens) {
    $sub = $d.Range($pos, $pos + $len)
    $sub.FormattedText = $sub.FormattedText
    $pos = $pos + $len
}
